$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Extend the data table with four new year columns (2020-2023): X:AA
#    Copy the formatting of the last existing year block (T:W) into the
#    new X:AA block for rows 4-7, then overwrite the copied values with
#    the real figures for the new years.
# ---------------------------------------------------------------------

$ws.Range("T4:W4").Copy($ws.Range("X4:AA4"))
$ws.Range("T5:W5").Copy($ws.Range("X5:AA5"))
$ws.Range("T6:W6").Copy($ws.Range("X6:AA6"))
$ws.Range("T7:W7").Copy($ws.Range("X7:AA7"))

# Row 4 - years
$ws.Range("X4").Value = 2020
$ws.Range("Y4").Value = 2021
$ws.Range("Z4").Value = 2022
$ws.Range("AA4").Value = 2023

# Row 5 - Investment loan
$ws.Range("X5").Value = 23780
$ws.Range("Y5").Value = 44660
$ws.Range("Z5").Value = 25000
$ws.Range("AA5").Value = 13010

# Row 6 - Investment grant
$ws.Range("X6").Value = 38240
$ws.Range("Y6").Value = 7950
$ws.Range("Z6").Value = 23000
$ws.Range("AA6").Value = 16390

# Row 7 - Investment loan and grant
$ws.Range("X7").Value = 62020
$ws.Range("Y7").Value = 52610
$ws.Range("Z7").Value = 48000
$ws.Range("AA7").Value = 29400

# ---------------------------------------------------------------------
# 2) The first set of year columns (D:K, i.e. 2000-2007) is no longer the
#    focus of the sheet now that newer years were appended - hide them,
#    matching the "hidden narrow columns" treatment applied upstream.
# ---------------------------------------------------------------------

$ws.Columns("D:K").Hidden = $true

# ---------------------------------------------------------------------
# 3) Bump the row heights of the header/data rows slightly so the taller
#    (now four-years-wider) table reads comfortably.
# ---------------------------------------------------------------------

$ws.Rows("4:4").RowHeight = 16.5
$ws.Rows("5:5").RowHeight = 16.5
$ws.Rows("6:6").RowHeight = 16.5
$ws.Rows("7:7").RowHeight = 16.5
